$d = $word.ActiveDocument
$rng = $d.Content
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>Card</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>：</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>接订单</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>Conversation</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>：</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>已注册用户，在赚外快页面接单，以便实现服务的供求匹配。</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>Confirmation</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:sz w:val="24"/></w:rPr><w:t>：</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>游客接单（不通过）</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>注册用户接单（通过）</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>点击订单显示详细信息（通过）</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>接单成功有提示（通过）</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>取快递订单在被接受前，部分信息隐藏（通过）</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>接单后，其他用户不能再次接单（通过）</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>接单后，下单用户收到反馈信息（通过）</w:t></w:r></w:p><w:p><w:r><w:t>PLACEHOLDER_TO_BE_CLEARED</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
Write-Host "ParaCount after insert:" $d.Paragraphs.Count

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRng = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$delRng.Delete()
Write-Host "ParaCount after cleanup:" $d.Paragraphs.Count
